$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 10, pushing existing rows 10-39 down to 12-41
$ws.Rows("10:11").Insert()

# Row 10 - new weekly price entry (Primera)
$ws.Range("A10").Value = 1
$ws.Range("B10").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C10").Value = "Arica y Parinacota"
$ws.Range("D10").Value = 44414
$ws.Range("E10").Value = 15
$ws.Range("F10").Value = 100112036
$ws.Range("G10").Value = "Caigua"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 130
$ws.Range("K10").Value = 8000
$ws.Range("L10").Value = 9000
$ws.Range("M10").Value = 8500
$ws.Range("N10").Value = "$/caja 20 kilos"
$ws.Range("O10").Value = "Región de Arica y Parinacota"
$ws.Range("P10").Value = 425
$ws.Range("Q10").Value = 20
$ws.Range("R10").Value = "Hortaliza"

# Row 11 - new weekly price entry (Segunda)
$ws.Range("A11").Value = 1
$ws.Range("B11").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C11").Value = "Arica y Parinacota"
$ws.Range("D11").Value = 44414
$ws.Range("E11").Value = 15
$ws.Range("F11").Value = 100112036
$ws.Range("G11").Value = "Caigua"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Segunda"
$ws.Range("J11").Value = 120
$ws.Range("K11").Value = 6000
$ws.Range("L11").Value = 7000
$ws.Range("M11").Value = 6500
$ws.Range("N11").Value = "$/caja 20 kilos"
$ws.Range("O11").Value = "Región de Arica y Parinacota"
$ws.Range("P11").Value = 325
$ws.Range("Q11").Value = 20
$ws.Range("R11").Value = "Hortaliza"
